$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (spaces -> underscores) for the two "Storage Capacity" columns
$ws.Range("G1").Value = "Storage_Capacity_Allocated_GB"
$ws.Range("H1").Value = "Storage_Capacity_Used_GB"

# Row 2: VMAX storage array (clear the stale E2:G2 values left over from the old row)
$ws.Range("A2").Value = 192602203
$ws.Range("B2").Value = "VMAX"
$ws.Range("C2").Value = "EMC"
$ws.Range("D2").Value = "VMAX-1"
$ws.Range("E2:G2").ClearContents()
$ws.Range("H2").Value = 71.3
$ws.Range("I2:P2").Value = "NULL"

# Row 3: Clariion CX3-40
$ws.Range("A3").Value = "APM00031901577"
$ws.Range("B3").Value = "Clariion"
$ws.Range("C3").Value = "EMC"
$ws.Range("D3").Value = "CX3-40"
$ws.Range("H3").Value = 12.4
$ws.Range("I3:P3").Value = "NULL"

# Row 4: Clariion CX3-80
$ws.Range("A4").Value = "APM00081100918"
$ws.Range("B4").Value = "Clariion"
$ws.Range("C4").Value = "EMC"
$ws.Range("D4").Value = "CX3-80"
$ws.Range("H4").Value = 36.2
$ws.Range("I4:P4").Value = "NULL"

# Row 5: Clariion CX4-240
$ws.Range("A5").Value = "APM00094701344"
$ws.Range("B5").Value = "Clariion"
$ws.Range("C5").Value = "EMC"
$ws.Range("D5").Value = "CX4-240"
$ws.Range("H5").Value = 22.5
$ws.Range("I5:P5").Value = "NULL"

# Row 6: Clariion CX4-960
$ws.Range("A6").Value = "APM00095103614"
$ws.Range("B6").Value = "Clariion"
$ws.Range("C6").Value = "EMC"
$ws.Range("D6").Value = "CX4-960"
$ws.Range("H6").Value = 228.4
$ws.Range("I6:P6").Value = "NULL"

# Row 7: Clariion CX4-480
$ws.Range("A7").Value = "APM00104904015"
$ws.Range("B7").Value = "Clariion"
$ws.Range("C7").Value = "EMC"
$ws.Range("D7").Value = "CX4-480"
$ws.Range("H7").Value = 103.7
$ws.Range("I7:P7").Value = "NULL"

# Row 8: VNX array, Model is numeric 7500
$ws.Range("A8").Value = "APM00114801908"
$ws.Range("B8").Value = "VNX"
$ws.Range("C8").Value = "EMC"
$ws.Range("D8").Value = 7500
$ws.Range("H8").Value = 169.6
$ws.Range("I8:P8").Value = "NULL"

# Selection moves to cover rows 2-8 (entire rows), matching the saved view state
$ws.Range("A2:A8").EntireRow.Select()
